$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the header style (bold, bordered, centered) to column A rows 2-7
# (matches the header's formatting, as seen in the target sheet) before
# writing values, so the paste-special doesn't clobber the new content.
$ws.Range("A1").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

# Header row
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# Data rows: eng/fra title master data
$data = @(
    @("eng", "MIR", "Mr",           "Male Title",                $true),
    @("eng", "MRS", "Mrs",          "Female Title",               $true),
    @("eng", "MIS", "Miss",         "Unmarried Female Title",     $true),
    @("fra", "MIR", "Monsieur",     "Titre masculin",             $true),
    @("fra", "MRS", "Madame",       "Titre féminin",              $true),
    @("fra", "MIS", "Mademoiselle", "Titre de femme célibataire", $true)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}

Write-Output $ws.Range("A1:E7").Value
